# SIT_TESTDATA.xlsx / LE service specification factory excel file
#
# Change: on the "LE_FTSP" sheet, test case row 3
# ("TC_002_Send the Invalid request") has its RunMode flag (column B)
# switched from "Yes" to "No" so that test case no longer runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LE_FTSP")

# RunMode for TC_002_Send the Invalid request (row 3) : Yes -> No
$ws.Range("B3").Value2 = "No"

# Reflect where the editor ended up after making the change (matches the
# saved cursor position in the source workbook).
$ws.Activate()
$ws.Range("B7").Select()
